$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('The Explorer', ['Hero', 'You may play an additional land on each of your turns.'])"
$ws.Range("A3").Value = "('The General', ['Hero', 'Exile The General: Creatures you control get +1/+1 until end of turn. Untap them.'])"
$ws.Range("A4").Value = "('The Provider', ['Hero', 'Exile The Provider: Put two +1/+1 counters on target creature you control. You gain life equal to that creature’s toughness.'])"
$ws.Range("A5").Value = "('The Savant', ['Hero', 'Exile The Savant: Tap all creatures your opponents control. Those creatures don’t untap during their controllers’ next untap steps.'])"
$ws.Range("A6").Value = "('The Tyrant', ['Hero', 'Exile The Tyrant: Creatures your opponents control get -1/-1 until end of turn.'])"
$ws.Range("A7").Value = "('The Vanquisher', ['Hero', 'Your starting hand size is increased by one.', 'Your maximum hand size is increased by one.'])"
$ws.Range("A8").Value = "('The Warmonger', ['Hero', 'Exile The Warmonger: Creatures you control get +2/+0 and gain haste until end of turn.'])"

$ws.Range("A9:A23").EntireRow.Delete()
